$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Cardiology" to "Session"
$ws.Name = "Session"

# Delete the two obsolete log rows (Student ID 211631 logged at 15:13:46,
# and Student ID 211217 logged at 15:14:21). Delete the later row first so
# that the row numbers of earlier rows stay valid.
$ws.Rows.Item(35).Delete()
$ws.Rows.Item(31).Delete()
